$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country/province stats (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 667572
$ws.Range("C4").Value = 19424
$ws.Range("E4").Value = 576480
$ws.Range("G4").Value = 1315
$ws.Range("H4").Value = 33903

# Row 7 - Francia
$ws.Range("D7").Value = 32812
$ws.Range("E7").Value = 97131
$ws.Range("F7").Value = 6248
$ws.Range("G7").Value = 753
$ws.Range("H7").Value = 17920

# Row 18 - Suiza
$ws.Range("D18").Value = 15900
$ws.Range("E18").Value = 9551

# Row 22 - Israel
$ws.Range("B22").Value = 12758
$ws.Range("C22").Value = 257
$ws.Range("D22").Value = 2818
$ws.Range("E22").Value = 9798
$ws.Range("F22").Value = 181
$ws.Range("G22").Value = 12
$ws.Range("H22").Value = 142

# Row 25 - Peru
$ws.Range("B25").Value = 12491
$ws.Range("C25").Value = 1016
$ws.Range("D25").Value = 6120
$ws.Range("E25").Value = 6097
$ws.Range("F25").Value = 169
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 274

# Row 50 - Luxemburgo
$ws.Range("D50").Value = 552
$ws.Range("E50").Value = 2823

# Rows 69/70 - Nueva Zelanda / Kazajistan swap order (Kazajistan now ranks
# above Nueva Zelanda) with updated figures
$ws.Range("A69").Value = "Kazajistan"
$ws.Range("B69").Value = 1402
$ws.Range("C69").Value = 107
$ws.Range("D69").Value = 277
$ws.Range("E69").Value = 1108
$ws.Range("F69").Value = 22
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 17

$ws.Range("A70").Value = "Nueva Zelanda"
$ws.Range("B70").Value = 1401
$ws.Range("C70").Value = 15
$ws.Range("D70").Value = 770
$ws.Range("E70").Value = 622
$ws.Range("F70").Value = 3
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 9

# Row 114 - Isla de Man
$ws.Range("B114").Value = 284
$ws.Range("C114").Value = 28
$ws.Range("D114").Value = 154

# Row 128 - Ruanda
$ws.Range("B128").Value = 138
$ws.Range("C128").Value = 2
$ws.Range("D128").Value = 60
$ws.Range("E128").Value = 78
